$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit re-shuffles the per-row price/date/quality data (columns D, I, J, K, L, M, P)
# across rows (weekly re-aggregation). Snapshot the original values first (using
# Value2, which reliably returns the underlying scalar rather than a property handle)
# before overwriting any cells, since source and destination rows overlap.
$D2 = $ws.Cells.Item(2, 4).Value2
$I2 = $ws.Cells.Item(2, 9).Value2
$J2 = $ws.Cells.Item(2, 10).Value2
$K2 = $ws.Cells.Item(2, 11).Value2
$L2 = $ws.Cells.Item(2, 12).Value2
$M2 = $ws.Cells.Item(2, 13).Value2
$P2 = $ws.Cells.Item(2, 16).Value2
$D3 = $ws.Cells.Item(3, 4).Value2
$I3 = $ws.Cells.Item(3, 9).Value2
$J3 = $ws.Cells.Item(3, 10).Value2
$K3 = $ws.Cells.Item(3, 11).Value2
$L3 = $ws.Cells.Item(3, 12).Value2
$M3 = $ws.Cells.Item(3, 13).Value2
$P3 = $ws.Cells.Item(3, 16).Value2
$D4 = $ws.Cells.Item(4, 4).Value2
$I4 = $ws.Cells.Item(4, 9).Value2
$J4 = $ws.Cells.Item(4, 10).Value2
$K4 = $ws.Cells.Item(4, 11).Value2
$L4 = $ws.Cells.Item(4, 12).Value2
$M4 = $ws.Cells.Item(4, 13).Value2
$P4 = $ws.Cells.Item(4, 16).Value2
$D5 = $ws.Cells.Item(5, 4).Value2
$I5 = $ws.Cells.Item(5, 9).Value2
$J5 = $ws.Cells.Item(5, 10).Value2
$K5 = $ws.Cells.Item(5, 11).Value2
$L5 = $ws.Cells.Item(5, 12).Value2
$M5 = $ws.Cells.Item(5, 13).Value2
$P5 = $ws.Cells.Item(5, 16).Value2
$D7 = $ws.Cells.Item(7, 4).Value2
$I7 = $ws.Cells.Item(7, 9).Value2
$J7 = $ws.Cells.Item(7, 10).Value2
$K7 = $ws.Cells.Item(7, 11).Value2
$L7 = $ws.Cells.Item(7, 12).Value2
$M7 = $ws.Cells.Item(7, 13).Value2
$P7 = $ws.Cells.Item(7, 16).Value2
$D8 = $ws.Cells.Item(8, 4).Value2
$I8 = $ws.Cells.Item(8, 9).Value2
$J8 = $ws.Cells.Item(8, 10).Value2
$K8 = $ws.Cells.Item(8, 11).Value2
$L8 = $ws.Cells.Item(8, 12).Value2
$M8 = $ws.Cells.Item(8, 13).Value2
$P8 = $ws.Cells.Item(8, 16).Value2
$D9 = $ws.Cells.Item(9, 4).Value2
$I9 = $ws.Cells.Item(9, 9).Value2
$J9 = $ws.Cells.Item(9, 10).Value2
$K9 = $ws.Cells.Item(9, 11).Value2
$L9 = $ws.Cells.Item(9, 12).Value2
$M9 = $ws.Cells.Item(9, 13).Value2
$P9 = $ws.Cells.Item(9, 16).Value2
$D10 = $ws.Cells.Item(10, 4).Value2
$I10 = $ws.Cells.Item(10, 9).Value2
$J10 = $ws.Cells.Item(10, 10).Value2
$K10 = $ws.Cells.Item(10, 11).Value2
$L10 = $ws.Cells.Item(10, 12).Value2
$M10 = $ws.Cells.Item(10, 13).Value2
$P10 = $ws.Cells.Item(10, 16).Value2
$D11 = $ws.Cells.Item(11, 4).Value2
$I11 = $ws.Cells.Item(11, 9).Value2
$J11 = $ws.Cells.Item(11, 10).Value2
$K11 = $ws.Cells.Item(11, 11).Value2
$L11 = $ws.Cells.Item(11, 12).Value2
$M11 = $ws.Cells.Item(11, 13).Value2
$P11 = $ws.Cells.Item(11, 16).Value2
$D12 = $ws.Cells.Item(12, 4).Value2
$I12 = $ws.Cells.Item(12, 9).Value2
$J12 = $ws.Cells.Item(12, 10).Value2
$K12 = $ws.Cells.Item(12, 11).Value2
$L12 = $ws.Cells.Item(12, 12).Value2
$M12 = $ws.Cells.Item(12, 13).Value2
$P12 = $ws.Cells.Item(12, 16).Value2
$D13 = $ws.Cells.Item(13, 4).Value2
$I13 = $ws.Cells.Item(13, 9).Value2
$J13 = $ws.Cells.Item(13, 10).Value2
$K13 = $ws.Cells.Item(13, 11).Value2
$L13 = $ws.Cells.Item(13, 12).Value2
$M13 = $ws.Cells.Item(13, 13).Value2
$P13 = $ws.Cells.Item(13, 16).Value2
$D14 = $ws.Cells.Item(14, 4).Value2
$I14 = $ws.Cells.Item(14, 9).Value2
$J14 = $ws.Cells.Item(14, 10).Value2
$K14 = $ws.Cells.Item(14, 11).Value2
$L14 = $ws.Cells.Item(14, 12).Value2
$M14 = $ws.Cells.Item(14, 13).Value2
$P14 = $ws.Cells.Item(14, 16).Value2
$D15 = $ws.Cells.Item(15, 4).Value2
$I15 = $ws.Cells.Item(15, 9).Value2
$J15 = $ws.Cells.Item(15, 10).Value2
$K15 = $ws.Cells.Item(15, 11).Value2
$L15 = $ws.Cells.Item(15, 12).Value2
$M15 = $ws.Cells.Item(15, 13).Value2
$P15 = $ws.Cells.Item(15, 16).Value2
$D18 = $ws.Cells.Item(18, 4).Value2
$I18 = $ws.Cells.Item(18, 9).Value2
$J18 = $ws.Cells.Item(18, 10).Value2
$K18 = $ws.Cells.Item(18, 11).Value2
$L18 = $ws.Cells.Item(18, 12).Value2
$M18 = $ws.Cells.Item(18, 13).Value2
$P18 = $ws.Cells.Item(18, 16).Value2
$D19 = $ws.Cells.Item(19, 4).Value2
$I19 = $ws.Cells.Item(19, 9).Value2
$J19 = $ws.Cells.Item(19, 10).Value2
$K19 = $ws.Cells.Item(19, 11).Value2
$L19 = $ws.Cells.Item(19, 12).Value2
$M19 = $ws.Cells.Item(19, 13).Value2
$P19 = $ws.Cells.Item(19, 16).Value2
$D20 = $ws.Cells.Item(20, 4).Value2
$I20 = $ws.Cells.Item(20, 9).Value2
$J20 = $ws.Cells.Item(20, 10).Value2
$K20 = $ws.Cells.Item(20, 11).Value2
$L20 = $ws.Cells.Item(20, 12).Value2
$M20 = $ws.Cells.Item(20, 13).Value2
$P20 = $ws.Cells.Item(20, 16).Value2
$D21 = $ws.Cells.Item(21, 4).Value2
$I21 = $ws.Cells.Item(21, 9).Value2
$J21 = $ws.Cells.Item(21, 10).Value2
$K21 = $ws.Cells.Item(21, 11).Value2
$L21 = $ws.Cells.Item(21, 12).Value2
$M21 = $ws.Cells.Item(21, 13).Value2
$P21 = $ws.Cells.Item(21, 16).Value2

# Apply the permuted values to their destination rows

# Row 2 gets the values previously held by row 3
$ws.Cells.Item(2, 4).Value = $D3
$ws.Cells.Item(2, 9).Value = $I3
$ws.Cells.Item(2, 10).Value = $J3
$ws.Cells.Item(2, 11).Value = $K3
$ws.Cells.Item(2, 12).Value = $L3
$ws.Cells.Item(2, 13).Value = $M3
$ws.Cells.Item(2, 16).Value = $P3

# Row 3 gets the values previously held by row 19
$ws.Cells.Item(3, 4).Value = $D19
$ws.Cells.Item(3, 9).Value = $I19
$ws.Cells.Item(3, 10).Value = $J19
$ws.Cells.Item(3, 11).Value = $K19
$ws.Cells.Item(3, 12).Value = $L19
$ws.Cells.Item(3, 13).Value = $M19
$ws.Cells.Item(3, 16).Value = $P19

# Row 4 gets the values previously held by row 20
$ws.Cells.Item(4, 4).Value = $D20
$ws.Cells.Item(4, 9).Value = $I20
$ws.Cells.Item(4, 10).Value = $J20
$ws.Cells.Item(4, 11).Value = $K20
$ws.Cells.Item(4, 12).Value = $L20
$ws.Cells.Item(4, 13).Value = $M20
$ws.Cells.Item(4, 16).Value = $P20

# Row 5 gets the values previously held by row 9
$ws.Cells.Item(5, 4).Value = $D9
$ws.Cells.Item(5, 9).Value = $I9
$ws.Cells.Item(5, 10).Value = $J9
$ws.Cells.Item(5, 11).Value = $K9
$ws.Cells.Item(5, 12).Value = $L9
$ws.Cells.Item(5, 13).Value = $M9
$ws.Cells.Item(5, 16).Value = $P9

# Row 7 gets the values previously held by row 14
$ws.Cells.Item(7, 4).Value = $D14
$ws.Cells.Item(7, 9).Value = $I14
$ws.Cells.Item(7, 10).Value = $J14
$ws.Cells.Item(7, 11).Value = $K14
$ws.Cells.Item(7, 12).Value = $L14
$ws.Cells.Item(7, 13).Value = $M14
$ws.Cells.Item(7, 16).Value = $P14

# Row 8 gets the values previously held by row 15
$ws.Cells.Item(8, 4).Value = $D15
$ws.Cells.Item(8, 9).Value = $I15
$ws.Cells.Item(8, 10).Value = $J15
$ws.Cells.Item(8, 11).Value = $K15
$ws.Cells.Item(8, 12).Value = $L15
$ws.Cells.Item(8, 13).Value = $M15
$ws.Cells.Item(8, 16).Value = $P15

# Row 9 gets the values previously held by row 7
$ws.Cells.Item(9, 4).Value = $D7
$ws.Cells.Item(9, 9).Value = $I7
$ws.Cells.Item(9, 10).Value = $J7
$ws.Cells.Item(9, 11).Value = $K7
$ws.Cells.Item(9, 12).Value = $L7
$ws.Cells.Item(9, 13).Value = $M7
$ws.Cells.Item(9, 16).Value = $P7

# Row 10 gets the values previously held by row 8
$ws.Cells.Item(10, 4).Value = $D8
$ws.Cells.Item(10, 9).Value = $I8
$ws.Cells.Item(10, 10).Value = $J8
$ws.Cells.Item(10, 11).Value = $K8
$ws.Cells.Item(10, 12).Value = $L8
$ws.Cells.Item(10, 13).Value = $M8
$ws.Cells.Item(10, 16).Value = $P8

# Row 11 gets the values previously held by row 5
$ws.Cells.Item(11, 4).Value = $D5
$ws.Cells.Item(11, 9).Value = $I5
$ws.Cells.Item(11, 10).Value = $J5
$ws.Cells.Item(11, 11).Value = $K5
$ws.Cells.Item(11, 12).Value = $L5
$ws.Cells.Item(11, 13).Value = $M5
$ws.Cells.Item(11, 16).Value = $P5

# Row 12 gets the values previously held by row 13
$ws.Cells.Item(12, 4).Value = $D13
$ws.Cells.Item(12, 9).Value = $I13
$ws.Cells.Item(12, 10).Value = $J13
$ws.Cells.Item(12, 11).Value = $K13
$ws.Cells.Item(12, 12).Value = $L13
$ws.Cells.Item(12, 13).Value = $M13
$ws.Cells.Item(12, 16).Value = $P13

# Row 13 gets the values previously held by row 10
$ws.Cells.Item(13, 4).Value = $D10
$ws.Cells.Item(13, 9).Value = $I10
$ws.Cells.Item(13, 10).Value = $J10
$ws.Cells.Item(13, 11).Value = $K10
$ws.Cells.Item(13, 12).Value = $L10
$ws.Cells.Item(13, 13).Value = $M10
$ws.Cells.Item(13, 16).Value = $P10

# Row 14 gets the values previously held by row 18
$ws.Cells.Item(14, 4).Value = $D18
$ws.Cells.Item(14, 9).Value = $I18
$ws.Cells.Item(14, 10).Value = $J18
$ws.Cells.Item(14, 11).Value = $K18
$ws.Cells.Item(14, 12).Value = $L18
$ws.Cells.Item(14, 13).Value = $M18
$ws.Cells.Item(14, 16).Value = $P18

# Row 15 gets the values previously held by row 21
$ws.Cells.Item(15, 4).Value = $D21
$ws.Cells.Item(15, 9).Value = $I21
$ws.Cells.Item(15, 10).Value = $J21
$ws.Cells.Item(15, 11).Value = $K21
$ws.Cells.Item(15, 12).Value = $L21
$ws.Cells.Item(15, 13).Value = $M21
$ws.Cells.Item(15, 16).Value = $P21

# Row 18 gets the values previously held by row 12
$ws.Cells.Item(18, 4).Value = $D12
$ws.Cells.Item(18, 9).Value = $I12
$ws.Cells.Item(18, 10).Value = $J12
$ws.Cells.Item(18, 11).Value = $K12
$ws.Cells.Item(18, 12).Value = $L12
$ws.Cells.Item(18, 13).Value = $M12
$ws.Cells.Item(18, 16).Value = $P12

# Row 19 gets the values previously held by row 4
$ws.Cells.Item(19, 4).Value = $D4
$ws.Cells.Item(19, 9).Value = $I4
$ws.Cells.Item(19, 10).Value = $J4
$ws.Cells.Item(19, 11).Value = $K4
$ws.Cells.Item(19, 12).Value = $L4
$ws.Cells.Item(19, 13).Value = $M4
$ws.Cells.Item(19, 16).Value = $P4

# Row 20 gets the values previously held by row 11
$ws.Cells.Item(20, 4).Value = $D11
$ws.Cells.Item(20, 9).Value = $I11
$ws.Cells.Item(20, 10).Value = $J11
$ws.Cells.Item(20, 11).Value = $K11
$ws.Cells.Item(20, 12).Value = $L11
$ws.Cells.Item(20, 13).Value = $M11
$ws.Cells.Item(20, 16).Value = $P11

# Row 21 gets the values previously held by row 2
$ws.Cells.Item(21, 4).Value = $D2
$ws.Cells.Item(21, 9).Value = $I2
$ws.Cells.Item(21, 10).Value = $J2
$ws.Cells.Item(21, 11).Value = $K2
$ws.Cells.Item(21, 12).Value = $L2
$ws.Cells.Item(21, 13).Value = $M2
$ws.Cells.Item(21, 16).Value = $P2
